$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.654473543167114
$ws.Range("B1").Value = 2.136698007583618
$ws.Range("C1").Value = 3.056734800338745
$ws.Range("D1").Value = 6.219033241271973
$ws.Range("E1").Value = 2.274832725524902
